$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("CO")
$ws2 = $wb.Worksheets.Item("SOUTIEN")

# Fill in the SOUTIEN sheet's still-open TBD slots (top portion first)
$ws2.Range("C2").Value = "William Baril"
$ws2.Range("C3").Value = "Clara Barbès"
$ws2.Range("C4").Value = "Stéphan Larose <br/>Mark Tremblay"
$ws2.Range("C6").Value = "Maurice Normand"

# Update the CO sheet's Cafeteria and Tour de la Releve rows
$ws1.Range("C9").Value = "France Galarneau<br/>Pierre Galarneau"
$ws1.Range("C23").Value = "Stephan Larose<br/>Marc Tremblay<br/>Stéphane Fiset<br/>Francis Fortin"

# Back to SOUTIEN for the remaining TBD slots
$ws2.Range("C7").Value = "Xavier Charron"
$ws2.Range("C9").Value = "Dessercom"

$ws2.Range("C11").Select()

$ws1.Activate()
$ws1.Range("C24").Select()
